$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 343, shifting existing rows 343:410 down to 344:411
$ws.Rows.Item(343).Insert()

# Populate the newly inserted row 343 with the new record's data
$ws.Cells.Item(343, 1).Value = 6
$ws.Cells.Item(343, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(343, 3).Value = "Metropolitana"
$ws.Cells.Item(343, 4).Value = 44644
$ws.Cells.Item(343, 5).Value = 13
$ws.Cells.Item(343, 6).Value = 100112043
$ws.Cells.Item(343, 7).Value = "Pepino ensalada"
$ws.Cells.Item(343, 8).Value = "Sin especificar"
$ws.Cells.Item(343, 9).Value = "Primera"
$ws.Cells.Item(343, 10).Value = 200
$ws.Cells.Item(343, 11).Value = 15000
$ws.Cells.Item(343, 12).Value = 17000
$ws.Cells.Item(343, 13).Value = 15800
$ws.Cells.Item(343, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(343, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(343, 16).Value = 263
$ws.Cells.Item(343, 17).Value = 60
$ws.Cells.Item(343, 18).Value = "Hortaliza"
